# "updated results from slr"
#
# Inserts 5 newly-reviewed papers at the top of the still-to-be-processed
# list in column A (rows 32-36), pushes the previously-listed papers down
# by 5 rows, and re-colors the boundary rows to flag progress:
#   - the first of the 5 new entries (row 32) gets a red highlight
#   - the row that used to be first in the old list ("UML-driven automated
#     software deployment", now row 37) gets a blue highlight

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the 5 new papers; this shifts rows 32:52 down to 37:57.
$ws.Rows("32:36").Insert()

# Excel's row insert propagates the formatting/borders of the row above
# into the freshly inserted row across its full used width - strip that
# back down to bare cells so only column A ends up populated, matching a
# simple list entry like the rest of the sheet.
$ws.Range("B32:E36").Clear()
$ws.Range("A32:A36").ClearFormats()

# New papers, entered in the order the author actually typed them (this
# is also the order they were interned into the shared-string table).
$ws.Range("A32").Value = "A Model-Driven Approach for Systematic Reproducibility and Replicability of Data Science Projects"
$ws.Range("A34").Value = "Industrial requirements for supporting AI-enhanced model-driven engineering"
$ws.Range("A35").Value = "MDE for machine learning-enabled software systems: a case study and comparison of MontiAnna &amp; ML-Quadrat"
$ws.Range("A36").Value = "Model-based fleet deployment in the IoT–edge–cloud continuum"
$ws.Range("A33").Value = "AI-augmented Model-Based Capabilities in the AIDOaRt Project: Continuous Development of Cyber-Physical Systems"

# Flag the first new row in red and the old top-of-list row (now pushed
# down to row 37) in blue.
$ws.Range("A32").Interior.Color = 255
$ws.Range("A37").Interior.Color = 12611584

# Leave the cursor where the author left it.
[void]$ws.Range("A11").Select()
